$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Text change in the installment paragraph:
#    " parcelas de R$ {{valor/num_parcelas}}"  ->  " parcelas de R$ {{valor_parcela}}"
#    (i.e. the placeholder "valor/num_parcelas" becomes "valor_parcela")
# ---------------------------------------------------------------------------
$content = $d.Content.Text
$oldSegment = " parcelas de R`$ {{valor/num_parcelas}}"
$idx = $content.IndexOf($oldSegment)
if ($idx -lt 0) {
    throw "Could not locate the installment text segment to replace."
}
$r = $d.Range($idx, $idx + $oldSegment.Length)
$r.Text = " parcelas de R`$ {{valor_parcela}}"

# ---------------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the signature-line paragraph to
#    the very end of the document (right after the final "{{cpf}}" merge
#    field, at the end of the CPF paragraph).
# ---------------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$pStart = $lastPara.Range.Start
$pEnd = $lastPara.Range.End - 1   # exclude the paragraph mark itself

$target = $d.Range($pStart, $pEnd)

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="007E5D97" w:rsidRPr="007E5D97" w:rsidRDefault="007E5D97" w:rsidP="008A3C3C"><w:pPr><w:ind w:firstLine="708"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/></w:rPr><w:t xml:space="preserve">CPF: </w:t></w:r><w:r w:rsidR="008A3C3C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="008A3C3C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/></w:rPr><w:t>cpf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="008A3C3C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/></w:rPr><w:t>}}</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xmlFrag)
